$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MarketCap")

# --- Insert 3 new rows before the old row 6 (pushes everything from the old
#     row 6 down to row 9, and everything below shifts down accordingly;
#     Excel auto-adjusts relative formula references on the shifted rows) ---
$ws.Rows("6:8").Insert()

# --- Updated market-cap figures (rows 2-4) ---
$ws.Range("C2").Value = 75163286
$ws.Range("F2").Value = "Index data"

$ws.Range("C3").Value = 43376587
$ws.Range("F3").Value = "Data as at: 31 Jun 2021"

$ws.Range("C4").Value = 31786699

# --- New rows 6-9: CRSP US Total Market Index cap split, with its own sum ---
# Row 6
$ws.Range("C2").Copy()
$ws.Range("C6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B6").Value = "CRSP US Total Market Index"
$ws.Range("C6").Value = 44358982

$ws.Range("D3").Copy()
$ws.Range("D6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D6").Formula = "=C6/C8"

# Row 7
$ws.Range("C2").Copy()
$ws.Range("C7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B7").Value = "FTSE Global All Cap ex US Index"
$ws.Range("C7").Value = 31786699

$ws.Range("D3").Copy()
$ws.Range("D7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D7").Formula = "=C7/C8"

# Row 8 - sum of the two new rows
$ws.Range("C2").Copy()
$ws.Range("C8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C8").Formula = "=SUM(C6:C7)"

$excel.CutCopyMode = 0

$wb.Save()
